$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '30.309.40'
Set-TextValue $ws.Range("E2") '  -0.29%  '

Set-TextValue $ws.Range("D3") '2.085.45'
Set-TextValue $ws.Range("E3") '  +3.41%  '

Set-TextValue $ws.Range("D4") '0.9999'
Set-TextValue $ws.Range("E4") '  -0.30%  '

Set-TextValue $ws.Range("D5") '329.01'
Set-TextValue $ws.Range("E5") '  +1.13%  '

Set-TextValue $ws.Range("D6") '0.9995'
Set-TextValue $ws.Range("E6") '  -0.19%  '

Set-TextValue $ws.Range("D7") '0.5231'
Set-TextValue $ws.Range("E7") '  +1.81%  '

Set-TextValue $ws.Range("D8") '0.4326'
Set-TextValue $ws.Range("E8") '  +2.67%  '

Set-TextValue $ws.Range("E9") '  +1.23%  '

Set-TextValue $ws.Range("D10") '46.82'
Set-TextValue $ws.Range("E10") '  +7.46%  '

Set-TextValue $ws.Range("E12") '  -1.14%  '

Set-TextValue $ws.Range("D13") '2.090.05'
Set-TextValue $ws.Range("E13") '  +3.82%  '

Set-TextValue $ws.Range("D14") '6.740'
Set-TextValue $ws.Range("E14") '  +2.13%  '

Set-TextValue $ws.Range("D15") '7.743'
Set-TextValue $ws.Range("E15") '  +3.59%  '

Set-TextValue $ws.Range("D16") '95.74'
Set-TextValue $ws.Range("E16") '  +1.17%  '

Set-TextValue $ws.Range("D17") '1.001'
Set-TextValue $ws.Range("E17") '  -0.17%  '

Set-TextValue $ws.Range("D18") '0.00001129'
Set-TextValue $ws.Range("E18") '  +1.32%  '

Set-TextValue $ws.Range("D19") '0.06639'
Set-TextValue $ws.Range("E19") '  +1.66%  '

Set-TextValue $ws.Range("D20") '18.93'
Set-TextValue $ws.Range("E20") '  -0.17%  '

Set-TextValue $ws.Range("D21") '0.9990'
Set-TextValue $ws.Range("E21") '  -0.21%  '

Set-TextValue $ws.Range("D22") '6.325'
Set-TextValue $ws.Range("E22") '  +1.90%  '

Set-TextValue $ws.Range("D23") '30.367.88'
Set-TextValue $ws.Range("E23") '  -0.28%  '

Set-TextValue $ws.Range("D24") '12.38'
Set-TextValue $ws.Range("E24") '  +4.47%  '

Set-TextValue $ws.Range("D25") '2.308'
Set-TextValue $ws.Range("E25") '  +3.01%  '

Set-TextValue $ws.Range("D26") '2.327.18'
Set-TextValue $ws.Range("E26") '  +3.42%  '

Set-TextValue $ws.Range("D27") '22.43'
Set-TextValue $ws.Range("E27") '  +0.10%  '

Set-TextValue $ws.Range("D28") '2.594'
Set-TextValue $ws.Range("E28") '  +6.74%  '

Set-TextValue $ws.Range("D29") '162.02'
Set-TextValue $ws.Range("E29") '  -0.44%  '

Set-TextValue $ws.Range("D30") '131.90'
Set-TextValue $ws.Range("E30") '  +0.43%  '

Set-TextValue $ws.Range("D31") '1.204'
Set-TextValue $ws.Range("E31") '  +5.67%  '

Set-TextValue $ws.Range("D32") '0.1071'
Set-TextValue $ws.Range("E32") '  +1.66%  '

Set-TextValue $ws.Range("D33") '1.668'
Set-TextValue $ws.Range("E33") '  +21.39%  '

Set-TextValue $ws.Range("D34") '6.188'
Set-TextValue $ws.Range("E34") '  +1.92%  '

Set-TextValue $ws.Range("D35") '3.876'
Set-TextValue $ws.Range("E35") '  +1.16%  '

Set-TextValue $ws.Range("D36") '9.917'
Set-TextValue $ws.Range("E36") '  +9.06%  '

Set-TextValue $ws.Range("D37") '0.02574'
Set-TextValue $ws.Range("E37") '  +1.55%  '

Set-TextValue $ws.Range("D38") '0.06683'
Set-TextValue $ws.Range("E38") '  +0.28%  '

Set-TextValue $ws.Range("B39") 'Aptos'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D39") '12.72'
Set-TextValue $ws.Range("E39") '  +3.02%  '

Set-TextValue $ws.Range("B40") 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D40") '5.455'
Set-TextValue $ws.Range("E40") '  -0.50%  '

Set-TextValue $ws.Range("E41") '  +2.92%  '

Set-TextValue $ws.Range("D42") '0.6834'
Set-TextValue $ws.Range("E42") '  +2.66%  '

Set-TextValue $ws.Range("D43") '1.243'
Set-TextValue $ws.Range("E43") '  +0.78%  '

Set-TextValue $ws.Range("D44") '0.9987'
Set-TextValue $ws.Range("E44") '  -0.21%  '

Set-TextValue $ws.Range("B45") 'Decentraland'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range("D45") '0.6393'
Set-TextValue $ws.Range("E45") '  +3.57%  '

Set-TextValue $ws.Range("B46") 'EnergySwap'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D46") '14.02'
Set-TextValue $ws.Range("E46") '  +2.64%  '

Set-TextValue $ws.Range("E47") '  +0.75%  '

Set-TextValue $ws.Range("D48") '3.611'
Set-TextValue $ws.Range("E48") '  -1.57%  '

Set-TextValue $ws.Range("D49") '1.251'
Set-TextValue $ws.Range("E49") '  -1.09%  '

Set-TextValue $ws.Range("E50") '  +7.43%  '

Set-TextValue $ws.Range("D51") '81.96'
Set-TextValue $ws.Range("E51") '  +1.16%  '
